$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.053.93"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.424.06"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "406.38"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.59"
$ws.Range("E6").Value = "  +3.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.594"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.691"
$ws.Range("E9").Value = "  +3.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("E10").Value = "  +4.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.98"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.44"
$ws.Range("E14").Value = "  -0.98%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.416.74"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.73"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.999.83"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("E19").Value = "  +11.61%  "
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "84.00"
$ws.Range("E21").Value = "  +1.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "313.87"
$ws.Range("E22").Value = "  +2.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.85"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.16"
$ws.Range("E24").Value = "  -1.90%  "
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.68"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.95"
$ws.Range("E27").Value = "  +6.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.15"
$ws.Range("E28").Value = "  -6.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.73"
$ws.Range("E29").Value = "  +6.14%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "43.79"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.116"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.35"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.68"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.318"
$ws.Range("E40").Value = "  +12.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "140.21"
$ws.Range("E41").Value = "  +4.22%  "
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.81"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.107.16"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.93"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.71"
$ws.Range("E51").Value = "  +18.17%  "
